$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 0.02258322285507441
$ws.Range("C2").Value = 0.004309184025731883
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 1.239511964969853

# Row 3 updates
$ws.Range("B3").Value = 0.001754667048134761
$ws.Range("C3").Value = 0.004309184025731883
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 0.6588563401661047
